$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1981132075471698
$ws.Range("C2").Value = 0.5471698113207547
$ws.Range("J2").Value = 0.02515723270440252
$ws.Range("P2").Value = 0.1446540880503145
$ws.Range("S2").Value = 0.08490566037735849
$ws.Range("B3").Value = 0.00558659217877095
$ws.Range("C3").Value = 0.01675977653631285
$ws.Range("J3").Value = 0.0223463687150838
$ws.Range("O3").Value = 0.00558659217877095
$ws.Range("P3").Value = 0.776536312849162
$ws.Range("S3").Value = 0.1731843575418995
$ws.Range("P4").Value = 0.675
$ws.Range("S4").Value = 0.325
$ws.Range("B6").Value = 0.03940886699507389
$ws.Range("D6").Value = 0.009852216748768473
$ws.Range("F6").Value = 0.04926108374384237
$ws.Range("J6").Value = 0.2266009852216749
$ws.Range("O6").Value = 0.01477832512315271
$ws.Range("Q6").Value = 0.1231527093596059
$ws.Range("R6").Value = 0.1182266009852217
$ws.Range("S6").Value = 0.4187192118226601
$ws.Range("B7").Value = 0.1408450704225352
$ws.Range("D7").Value = 0.01408450704225352
$ws.Range("F7").Value = 0.06338028169014084
$ws.Range("J7").Value = 0.1619718309859155
$ws.Range("O7").Value = 0.007042253521126761
$ws.Range("Q7").Value = 0.1619718309859155
$ws.Range("R7").Value = 0.1126760563380282
$ws.Range("S7").Value = 0.3380281690140845
$ws.Range("B8").Value = 0.09732360097323602
$ws.Range("D8").Value = 0.0121654501216545
$ws.Range("F8").Value = 0.06082725060827251
$ws.Range("J8").Value = 0.08759124087591241
$ws.Range("O8").Value = 0.0121654501216545
$ws.Range("Q8").Value = 0.2092457420924574
$ws.Range("R8").Value = 0.1386861313868613
$ws.Range("S8").Value = 0.3819951338199513
$ws.Range("B9").Value = 0.1690821256038647
$ws.Range("D9").Value = 0.004830917874396135
$ws.Range("F9").Value = 0.04830917874396135
$ws.Range("J9").Value = 0.0966183574879227
$ws.Range("O9").Value = 0.01932367149758454
$ws.Range("Q9").Value = 0.1642512077294686
$ws.Range("R9").Value = 0.106280193236715
$ws.Range("S9").Value = 0.391304347826087
$ws.Range("B10").Value = 0.1235584843492586
$ws.Range("D10").Value = 0.02471169686985173
$ws.Range("E10").Value = 0.001647446457990115
$ws.Range("F10").Value = 0.07660626029654036
$ws.Range("J10").Value = 0.09555189456342669
$ws.Range("O10").Value = 0.01153212520593081
$ws.Range("Q10").Value = 0.1927512355848435
$ws.Range("R10").Value = 0.1268533772652389
$ws.Range("S10").Value = 0.3467874794069193
$ws.Range("G11").Value = 0.187192118226601
$ws.Range("J11").Value = 0.05911330049261083
$ws.Range("K11").Value = 0.2167487684729064
$ws.Range("L11").Value = 0.5172413793103449
$ws.Range("S11").Value = 0.01970443349753695
$ws.Range("G12").Value = 0.7247706422018348
$ws.Range("J12").Value = 0.1834862385321101
$ws.Range("K12").Value = 0.01834862385321101
$ws.Range("L12").Value = 0.03669724770642202
$ws.Range("S12").Value = 0.03669724770642202
$ws.Range("G13").Value = 0.6744186046511628
$ws.Range("J13").Value = 0.2325581395348837
$ws.Range("S13").Value = 0.09302325581395349
$ws.Range("F15").Value = 0.01005025125628141
$ws.Range("H15").Value = 0.1608040201005025
$ws.Range("I15").Value = 0.1105527638190955
$ws.Range("J15").Value = 0.3768844221105528
$ws.Range("K15").Value = 0.03015075376884422
$ws.Range("M15").Value = 0.01005025125628141
$ws.Range("N15").Value = 0.01005025125628141
$ws.Range("O15").Value = 0.05527638190954774
$ws.Range("S15").Value = 0.2361809045226131
$ws.Range("F16").Value = 0.004830917874396135
$ws.Range("H16").Value = 0.1739130434782609
$ws.Range("I16").Value = 0.0966183574879227
$ws.Range("J16").Value = 0.3768115942028986
$ws.Range("K16").Value = 0.08695652173913043
$ws.Range("M16").Value = 0.02415458937198068
$ws.Range("O16").Value = 0.09178743961352658
$ws.Range("S16").Value = 0.1449275362318841
$ws.Range("F17").Value = 0.01272264631043257
$ws.Range("H17").Value = 0.1628498727735369
$ws.Range("I17").Value = 0.1221374045801527
$ws.Range("J17").Value = 0.4045801526717557
$ws.Range("K17").Value = 0.07633587786259542
$ws.Range("M17").Value = 0.01526717557251908
$ws.Range("O17").Value = 0.06615776081424936
$ws.Range("S17").Value = 0.1399491094147583
$ws.Range("F18").Value = 0.02962962962962963
$ws.Range("H18").Value = 0.1851851851851852
$ws.Range("I18").Value = 0.07407407407407407
$ws.Range("J18").Value = 0.4407407407407408
$ws.Range("M18").Value = 0.02222222222222222
$ws.Range("O18").Value = 0.07777777777777778
$ws.Range("S18").Value = 0.1037037037037037
$ws.Range("F19").Value = 0.008620689655172414
$ws.Range("H19").Value = 0.1982758620689655
$ws.Range("I19").Value = 0.08362068965517241
$ws.Range("J19").Value = 0.4318965517241379
$ws.Range("K19").Value = 0.07155172413793104
$ws.Range("M19").Value = 0.02327586206896552
$ws.Range("O19").Value = 0.06379310344827586
$ws.Range("S19").Value = 0.1189655172413793
